$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 31 (Informal economies and informal employment) - refreshed figures
# for "Percentage of informal employment in total employment" columns
# (C = total, D = females, E = males) on a handful of country / aggregate
# rows, following an update to the underlying source data.

# Row 50 - GMB / Gambia
$ws.Range("C50").Value = 81.599999999999994
$ws.Range("D50").Value = 86.9
$ws.Range("E50").Value = 76

# Row 61 - West Africa
$ws.Range("C61").Value = 87.65
$ws.Range("D61").Value = 90.985714285714295
$ws.Range("E61").Value = 84.8642857142857

# Row 62 - Africa
$ws.Range("C62").Value = 81.897560975609807
$ws.Range("D62").Value = 83.921951219512195
$ws.Range("E62").Value = 79.914634146341498

# Row 63 - World outside Africa
$ws.Range("C63").Value = 40.762790697674397
$ws.Range("D63").Value = 40.327906976744202
$ws.Range("E63").Value = 41.0162790697674

# Row 64 - Latin America and Caribbean
$ws.Range("C64").Value = 56.359090909090902
$ws.Range("D64").Value = 54.731818181818198
$ws.Range("E64").Value = 57.572727272727299

# Row 66 - World
$ws.Range("C66").Value = 54.042519685039402
$ws.Range("D66").Value = 54.4015748031496
$ws.Range("E66").Value = 53.574015748031499

# Row 68 - CEN-SAD
$ws.Range("C68").Value = 83.94
$ws.Range("D68").Value = 86.1
$ws.Range("E68").Value = 81.915000000000006

# Row 71 - ECOWAS
$ws.Range("C71").Value = 87.65
$ws.Range("D71").Value = 90.985714285714295
$ws.Range("E71").Value = 84.8642857142857

# Row 82 - Africa, Non-resource-rich countries
$ws.Range("C82").Value = 81.5513513513514
$ws.Range("D82").Value = 83.5324324324325
$ws.Range("E82").Value = 79.583783783783801

# Row 83 - ROW, Non-resource-rich countries
$ws.Range("C83").Value = 38.792307692307702
$ws.Range("D83").Value = 38.729487179487201
$ws.Range("E83").Value = 38.9

# Row 84 - Africa, Low income countries
$ws.Range("C84").Value = 91.0833333333333
$ws.Range("D84").Value = 93.938888888888897
$ws.Range("E84").Value = 88.605555555555597

# Row 90 - High income countries
$ws.Range("C90").Value = 14.0432432432432
$ws.Range("D90").Value = 13.4
$ws.Range("E90").Value = 14.524324324324301

# Row 91 - Africa, Least Developed Countries
$ws.Range("C91").Value = 89.040740740740802
$ws.Range("D91").Value = 92.248148148148204
$ws.Range("E91").Value = 86.551851851851893

# Row 97 - Africa, Fragile States
$ws.Range("C97").Value = 86.969230769230805
$ws.Range("D97").Value = 90.838461538461601
$ws.Range("E97").Value = 83.873076923076894
